# Implemented the config for Set create
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Insert a new row at position 3, shifting existing rows 3-10 down to 4-11
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).ClearFormats()

# Fill in the new row 3 values (order matters for shared-string table indices)
$ws.Range("B3").Value = "CreateSet_Test"
$ws.Range("A3").Value = "Create_Set"
$ws.Range("E3").Value = "Create Set1,Disease,kera"

# Match the row height used by similarly-formatted header row (row 2)
$ws.Rows.Item(3).RowHeight = 37.5

# Column B formatting: unhide and set to best-fit width instead of fixed hidden width
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(2).AutoFit()

# Expand the Table2 list object to include the newly inserted row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E11"))

# Update the selection to reflect where the user ended up (D3)
$ws.Range("D3").Select()
